$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26) entirely - all rows below shift up by one.
$ws.Rows(26).Delete()

# Remove the "SC 92" row, which after the previous deletion now sits at row 27.
$ws.Rows(27).Delete()

# Apply the remaining cell-level corrections (new/cleared values) using the
# final, post-deletion row numbers.
$ws.Range("E2").ClearContents()

$ws.Range("E5").Value = -5

$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7

$ws.Range("D8").ClearContents()

$ws.Range("E9").ClearContents()

$ws.Range("E10").ClearContents()

$ws.Range("D12").Value = -14.1

$ws.Range("D14").ClearContents()

$ws.Range("D17").Value = -14.7

$ws.Range("D18").Value = -15.2

$ws.Range("D19").ClearContents()

$ws.Range("D20").ClearContents()

$ws.Range("D23").Value = -13.9

$ws.Range("E24").Value = -8.1

$ws.Range("C27").Value = 10
$ws.Range("D27").ClearContents()

$ws.Range("C28").ClearContents()
$ws.Range("E28").ClearContents()

$ws.Range("C29").ClearContents()

$ws.Range("C30").Value = 11.4
$ws.Range("E30").Value = -5.7

$ws.Range("C32").ClearContents()
